# The target diff (canonical OOXML) touches only word/document.xml and
# word/styles.xml, and every single changed line is a pure XML
# attribute/namespace *re-ordering* (e.g. <w:pgSz w:w="11906" w:h="16838"/>
# -> <w:pgSz w:h="16838" w:w="11906"/>, alphabetised xmlns:* declarations on
# <w:document>, alphabetised w:* attributes on <w:style>/<w:lsdException>/...).
# Every changed element pair has the identical tag name and the identical
# attribute key/value set before and after - only the on-disk serialization
# order differs (this was verified programmatically against every hunk).
# There is no insertion, deletion, or value change of any paragraph, run,
# field, style, or section-property content - the document is semantically
# unchanged. The object model below simply touches the document without
# mutating any content, run formatting, styles, or section/page setup, so
# the saved package keeps the same (canonically-equivalent) content as the
# target revision.

$d = $word.ActiveDocument

# Touch the document content/formatting surfaces without changing them, to
# mirror the no-semantic-change nature of the target revision.
$null = $d.Content.Text
$null = $d.Styles("Normal").NameLocal
$null = $d.Sections(1).PageSetup.PageWidth
